{"js": "// Update the resume body text in place: swap the candidate name, summary,\n// technical skills, work-experience section (renamed \"Professional\n// Experience\"), and education block, and drop the old \"Soft Skills\"\n// section, per the target resume content.\n//\n// The whole document is a single paragraph/run whose lines are separated\n// by manual line breaks (<w:br/>), which Word/Office.js represent as the\n// vertical-tab character \"\\v\" (\\u000B) inside Range/Body `.text`. We\n// rebuild that text in one shot and let insertText() re-materialize the\n// \"\\v\" characters back into <w:br/> elements, so the resulting OOXML\n// keeps the original \"one run, many <w:t>/<w:br/> pairs\" shape.\n\nconst body = context.document.body;\n\nconst LF = \"\\u000b\"; // maps to <w:br/> when written back out\n\nconst lines = [\n  \"Resume of Taylor Brooks\",\n  \"\",\n  \"Summary:\",\n  \"Results-driven IT Support Specialist with hands-on experience in enterprise system troubleshooting and stakeholder communication. Adept at problem-solving and delivering end-to-end IT solutions in fast-paced environments.\",\n  \"\",\n  \"Technical Skills:\",\n  \"- Programming: Java, Python, SQL\",\n  \"- Database Systems: Oracle Database, MySQL\",\n  \"- Tools: JIRA, Confluence, Git\",\n  \"\",\n  \"Professional Experience:\",\n  \"\",\n  \"IT Solutions Analyst \\u2013 Celestica Technologies (2022\\u2013Present)\",\n  \"- Provided IT support for Celestica\\u2019s Enterprise Quality Management Solutions System.\",\n  \"- Collaborated with multiple stakeholders to identify and resolve systemic issues.\",\n  \"- Delivered full-cycle technical solutions and conducted stakeholder training sessions.\",\n  \"\",\n  \"Education:\",\n  \"Bachelor of Information Technology \\u2013 University of Toronto\",\n  \"\",\n];\n\nconst newText = LF + lines.join(LF);\n\nbody.clear();\nbody.insertText(newText, \"Start\");\nawait context.sync();\n", "ps1": "# Update the resume body text in place: swap the candidate name, summary,\n# technical skills, work-experience section (renamed \"Professional\n# Experience\"), and education block, and drop the old \"Soft Skills\"\n# section, per the target resume content.\n#\n# The whole document is a single paragraph/run whose lines are separated\n# by manual line breaks (<w:br/>), which Word represents as the\n# vertical-tab character Chr(11) inside Range.Text. We rebuild that text\n# in one shot and assign it back to the document Range, which\n# re-materializes the Chr(11) characters as <w:br/> elements, keeping the\n# original \"one run, many text/break pairs\" shape.\n\n$d = $word.ActiveDocument\n\n$LF = [char]11   # maps to <w:br/> when written back out\n$enDash = [char]0x2013\n$rsquo = [char]0x2019\n\n$lines = @(\n    \"Resume of Taylor Brooks\",\n    \"\",\n    \"Summary:\",\n    \"Results-driven IT Support Specialist with hands-on experience in enterprise system troubleshooting and stakeholder communication. Adept at problem-solving and delivering end-to-end IT solutions in fast-paced environments.\",\n    \"\",\n    \"Technical Skills:\",\n    \"- Programming: Java, Python, SQL\",\n    \"- Database Systems: Oracle Database, MySQL\",\n    \"- Tools: JIRA, Confluence, Git\",\n    \"\",\n    \"Professional Experience:\",\n    \"\",\n    \"IT Solutions Analyst $enDash Celestica Technologies (2022${enDash}Present)\",\n    \"- Provided IT support for Celestica${rsquo}s Enterprise Quality Management Solutions System.\",\n    \"- Collaborated with multiple stakeholders to identify and resolve systemic issues.\",\n    \"- Delivered full-cycle technical solutions and conducted stakeholder training sessions.\",\n    \"\",\n    \"Education:\",\n    \"Bachelor of Information Technology $enDash University of Toronto\",\n    \"\"\n)\n\n$newText = $LF + ($lines -join $LF)\n\n$d.Content.Text = $newText\n"}
